$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Foglio1")
$ws2 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1 ("Sheet1" tab): append a province -> district-count summary table ---

$ws2.Range("A9").Value = "Province"
$ws2.Range("B9").Value = "Number of districts"

$provinces = @(
    "Milano", "Bergamo", "Brescia", "Como", "Cremona", "Lodi e Crema",
    "Mantova", "Pavia", "Sondrio", "Venezia ", "Belluno", "Padova   ",
    "Rovigo", "Treviso", "Udine  ", "Verona ", "Vicenza     "
)

$row = 10
foreach ($prov in $provinces) {
    $ws2.Range("A$row").Value = $prov
    $ws2.Range("B$row").Formula = "=COUNTIFS(Foglio1!`$A`$2:`$A`$182,Sheet1!`$A$row)"
    $row = $row + 1
}

# --- Sheet1 ("Foglio1" tab): scroll / selection state change ---

$ws1.Range("A148").Select()
$ws1.Range("A1:A182").Select()

# --- Sheet2 ("Sheet1" tab): selection state change ---

$ws2.Range("F16").Select()

$wb.Save()
